$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated solar energy production predictions (model retrained to fix
# model-existence issue in Streamlit Cloud deployment).
$ws.Range("C12").Value = 0.188
$ws.Range("C13").Value = 0.5610000000000001
$ws.Range("C15").Value = 1.171
$ws.Range("C16").Value = 1.2
$ws.Range("C17").Value = 1.28
$ws.Range("C18").Value = 1.139
$ws.Range("C19").Value = 0.8129999999999999
$ws.Range("C20").Value = 0.25
$ws.Range("C21").Value = 0.041
$ws.Range("C36").Value = 0.15
$ws.Range("C37").Value = 0.586
$ws.Range("C38").Value = 1.326
$ws.Range("C39").Value = 1.181
$ws.Range("C40").Value = 1.379
$ws.Range("C41").Value = 1.404
$ws.Range("C42").Value = 1.115
$ws.Range("C43").Value = 0.753
$ws.Range("C44").Value = 0.245
$ws.Range("C63").Value = 1.572
$ws.Range("C64").Value = 2.125
$ws.Range("C65").Value = 2.006
$ws.Range("C66").Value = 1.307
$ws.Range("C84").Value = 0.068
$ws.Range("C85").Value = 0.419
$ws.Range("C86").Value = 0.891
$ws.Range("C87").Value = 1.152
$ws.Range("C88").Value = 1.304
$ws.Range("C90").Value = 0.87
$ws.Range("C91").Value = 0.521
$ws.Range("C92").Value = 0.165
$ws.Range("C109").Value = 0.196
$ws.Range("C110").Value = 0.41
$ws.Range("C111").Value = 0.646
$ws.Range("C112").Value = 0.593
$ws.Range("C113").Value = 0.43
$ws.Range("C114").Value = 0.177
$ws.Range("C115").Value = 0.126
$ws.Range("C116").Value = 0.055
$ws.Range("C117").Value = 0.012
$ws.Range("C133").Value = 0.5620000000000001
$ws.Range("C134").Value = 1.345
$ws.Range("C135").Value = 1.591
$ws.Range("C136").Value = 2.159
$ws.Range("C137").Value = 2.071
$ws.Range("C138").Value = 1.406
$ws.Range("C139").Value = 0.962
$ws.Range("C140").Value = 0.263
$ws.Range("C141").Value = 0.029
$ws.Range("C156").Value = 0.099
$ws.Range("C157").Value = 0.434
$ws.Range("C158").Value = 1.242
$ws.Range("C159").Value = 1.531
$ws.Range("C160").Value = 2.014
$ws.Range("C161").Value = 1.817
$ws.Range("C162").Value = 1.358
$ws.Range("C163").Value = 0.726
$ws.Range("C164").Value = 0.248
$ws.Range("C165").Value = 0.027
